$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores every cell in columns B:E as text (the "Price" and
# "Volume(1h)" columns use separators/signs/padding that only work as text).
# Some of the new Price values look like plain numbers though (e.g. "320.72"),
# so a normal Value assignment would be auto-typed as a number by Excel.
# Temporarily mark those cells as Text format before writing to them, then
# drop the format again afterwards so the cell style ends up unchanged.
$textCells = @(
  "D5", "D6", "D7", "D9", "D11", "D12", "D14", "D16", "D19", "D20", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D38", "D40", "D41", "D42", "D45", "D46", "D47", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '44.305.19'
$ws.Range("E2").Value = '  +3.79%  '
$ws.Range("D3").Value = '2.285.11'
$ws.Range("E3").Value = '  +3.43%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").Value = '320.72'
$ws.Range("E5").Value = '  +2.15%  '
$ws.Range("D6").Value = '104.48'
$ws.Range("E6").Value = '  +6.74%  '
$ws.Range("D7").Value = '0.591'
$ws.Range("E7").Value = '  +1.97%  '
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").Value = '0.575'
$ws.Range("E9").Value = '  +3.17%  '
$ws.Range("E10").Value = '  +7.42%  '
$ws.Range("D11").Value = '0.0843'
$ws.Range("E11").Value = '  +2.56%  '
$ws.Range("D12").Value = '7.93'
$ws.Range("E12").Value = '  +2.54%  '
$ws.Range("E13").Value = '  +2.40%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '0.884'
$ws.Range("E14").Value = '  +2.85%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.630.94'
$ws.Range("E15").Value = '  +3.36%  '
$ws.Range("D16").Value = '14.67'
$ws.Range("E16").Value = '  +4.26%  '
$ws.Range("D17").Value = '2.286.34'
$ws.Range("E17").Value = '  +3.97%  '
$ws.Range("D18").Value = '44.199.43'
$ws.Range("E18").Value = '  +3.81%  '
$ws.Range("D19").Value = '14.22'
$ws.Range("E19").Value = '  -3.48%  '
$ws.Range("D20").Value = '0.0000100'
$ws.Range("E20").Value = '  +4.83%  '
$ws.Range("E21").Value = '  +3.77%  '
$ws.Range("D22").Value = '66.49'
$ws.Range("E22").Value = '  +2.44%  '
$ws.Range("E23").Value = '  +3.13%  '
$ws.Range("D24").Value = '238.17'
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("D25").Value = '2.21'
$ws.Range("E25").Value = '  +4.58%  '
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").Value = '10.28'
$ws.Range("E27").Value = '  +2.78%  '
$ws.Range("D28").Value = '39.38'
$ws.Range("E28").Value = '  +16.33%  '
$ws.Range("D29").Value = '2.22'
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").Value = '6.55'
$ws.Range("E30").Value = '  +5.16%  '
$ws.Range("D31").Value = '163.73'
$ws.Range("E31").Value = '  +5.59%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.0888'
$ws.Range("E32").Value = '  +1.64%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '20.57'
$ws.Range("E33").Value = '  +1.13%  '
$ws.Range("E34").Value = '  -1.23%  '
$ws.Range("D35").Value = '2.10'
$ws.Range("E35").Value = '  +6.64%  '
$ws.Range("D36").Value = '3.30'
$ws.Range("E36").Value = '  +4.40%  '
$ws.Range("E37").Value = '  +11.63%  '
$ws.Range("D38").Value = '0.122'
$ws.Range("E38").Value = '  -0.24%  '
$ws.Range("E39").Value = '  +2.47%  '
$ws.Range("D40").Value = '3.95'
$ws.Range("E40").Value = '  +6.88%  '
$ws.Range("D41").Value = '15.72'
$ws.Range("E41").Value = '  +30.11%  '
$ws.Range("D42").Value = '0.0328'
$ws.Range("E42").Value = '  +1.35%  '
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").Value = '1.776.07'
$ws.Range("E44").Value = '  -4.97%  '
$ws.Range("D45").Value = '0.209'
$ws.Range("E45").Value = '  +1.34%  '
$ws.Range("D46").Value = '85.72'
$ws.Range("E46").Value = '  -3.46%  '
$ws.Range("D47").Value = '5.43'
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("D49").Value = '75.63'
$ws.Range("E49").Value = '  +0.30%  '
$ws.Range("D50").Value = '59.92'
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").Value = '104.95'
$ws.Range("E51").Value = '  +4.15%  '

foreach ($addr in $textCells) {
  $ws.Range($addr).Style = "Normal"
}
